$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 38.916668
$ws.Cells.Item(38, 9).Value = 38.916668
$ws.Cells.Item(38, 11).Value = 116.750004
$ws.Cells.Item(38, 13).Value = 255.249996
$ws.Cells.Item(39, 8).Value = 356.33334
$ws.Cells.Item(39, 9).Value = 429.57144
$ws.Cells.Item(39, 10).Value = 100
$ws.Cells.Item(39, 11).Value = 1288.71432
$ws.Cells.Item(39, 12).Value = 300
$ws.Cells.Item(39, 13).Value = -992.71432
$ws.Cells.Item(39, 14).Value = -892
$ws.Cells.Item(42, 8).Value = 465
$ws.Cells.Item(42, 9).Value = 202.66667
$ws.Cells.Item(42, 10).Value = 596.1667
$ws.Cells.Item(42, 11).Value = 608.00001
$ws.Cells.Item(42, 12).Value = 1788.5001
$ws.Cells.Item(42, 13).Value = -378.00001
$ws.Cells.Item(42, 14).Value = -2248.5001
$ws.Cells.Item(53, 8).Value = 282.72223
$ws.Cells.Item(53, 9).Value = 167.91667
$ws.Cells.Item(53, 11).Value = 167.91667
$ws.Cells.Item(53, 13).Value = 469.08333
$ws.Cells.Item(99, 8).Value = 1582
$ws.Cells.Item(99, 9).Value = 1582
$ws.Cells.Item(99, 11).Value = 4746
$ws.Cells.Item(99, 13).Value = -3248
$ws.Cells.Item(101, 8).Value = 2108.7144
$ws.Cells.Item(101, 9).Value = 2773.2
$ws.Cells.Item(101, 10).Value = 447.5
$ws.Cells.Item(101, 11).Value = 8319.599999999999
$ws.Cells.Item(101, 12).Value = 1342.5
$ws.Cells.Item(101, 13).Value = -6697.599999999999
$ws.Cells.Item(101, 14).Value = -4586.5
$ws.Cells.Item(108, 8).Value = 54999
$ws.Cells.Item(108, 10).Value = 54999
$ws.Cells.Item(108, 12).Value = 54999
$ws.Cells.Item(108, 14).Value = -62679
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).ClearContents()
$ws.Cells.Item(109, 14).Value = 0
$ws.Cells.Item(112, 8).Value = 2057.5625
$ws.Cells.Item(112, 10).Value = 2232.5386
$ws.Cells.Item(112, 12).Value = 6697.6158
$ws.Cells.Item(112, 14).Value = -8913.6158
$ws.Cells.Item(138, 8).Value = 3976.7334
$ws.Cells.Item(138, 9).Value = 3693.7144
$ws.Cells.Item(138, 10).Value = 4224.375
$ws.Cells.Item(138, 11).Value = 11081.1432
$ws.Cells.Item(138, 12).Value = 12673.125
$ws.Cells.Item(138, 13).Value = -5941.143199999999
$ws.Cells.Item(138, 14).Value = -22953.125

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(53, 8).Value = 10000
$ws.Cells.Item(53, 9).Value = 10000
$ws.Cells.Item(53, 11).Value = 10000
$ws.Cells.Item(53, 13).Value = -9318
$ws.Cells.Item(132, 8).Value = 2100
$ws.Cells.Item(132, 9).Value = 2100
$ws.Cells.Item(132, 11).Value = 6300
$ws.Cells.Item(132, 13).Value = -3770

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3829.92
$ws.Cells.Item(86, 9).Value = 3199.889
$ws.Cells.Item(86, 10).Value = 5450
$ws.Cells.Item(86, 11).Value = 3199.889
$ws.Cells.Item(86, 12).Value = 5450
$ws.Cells.Item(86, 13).Value = -2076.889
$ws.Cells.Item(86, 14).Value = -7696
$ws.Cells.Item(89, 8).Value = 3829.92
$ws.Cells.Item(89, 9).Value = 3199.889
$ws.Cells.Item(89, 10).Value = 5450
$ws.Cells.Item(89, 11).Value = 15999.445
$ws.Cells.Item(89, 12).Value = 27250
$ws.Cells.Item(89, 13).Value = -10383.445
$ws.Cells.Item(89, 14).Value = -38482
$ws.Cells.Item(105, 8).Value = 2261.3076
$ws.Cells.Item(105, 9).Value = 1824.75
$ws.Cells.Item(105, 11).Value = 1824.75
$ws.Cells.Item(105, 13).Value = -77.75
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).ClearContents()
$ws.Cells.Item(110, 14).Value = 0
$ws.Cells.Item(134, 8).Value = 2562.3635
$ws.Cells.Item(134, 9).Value = 2365.111
$ws.Cells.Item(134, 10).Value = 3450
$ws.Cells.Item(134, 11).Value = 7095.333
$ws.Cells.Item(134, 12).Value = 10350
$ws.Cells.Item(134, 13).Value = -4560.333
$ws.Cells.Item(134, 14).Value = -15420

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(8, 8).Value = 2163.3333
$ws.Cells.Item(8, 10).Value = 2163.3333
$ws.Cells.Item(8, 12).Value = 2163.3333
$ws.Cells.Item(8, 14).Value = -2443.3333
$ws.Cells.Item(22, 8).Value = 1986
$ws.Cells.Item(22, 9).Value = 847.25
$ws.Cells.Item(22, 10).Value = 3124.75
$ws.Cells.Item(22, 11).Value = 847.25
$ws.Cells.Item(22, 12).Value = 3124.75
$ws.Cells.Item(22, 13).Value = -497.25
$ws.Cells.Item(22, 14).Value = -3824.75
$ws.Cells.Item(31, 8).Value = 2302.348
$ws.Cells.Item(31, 9).Value = 2179.7273
$ws.Cells.Item(31, 11).Value = 2179.7273
$ws.Cells.Item(31, 13).Value = -1884.7273
$ws.Cells.Item(33, 8).Value = 2143.75
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(34, 8).Value = 2302.348
$ws.Cells.Item(34, 9).Value = 2179.7273
$ws.Cells.Item(34, 11).Value = 2179.7273
$ws.Cells.Item(34, 13).Value = -1977.7273

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 100.5
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 170
$ws.Cells.Item(51, 8).Value = 946.75
$ws.Cells.Item(51, 9).Value = 946.75
$ws.Cells.Item(51, 11).Value = 2840.25
$ws.Cells.Item(51, 13).Value = -2380.25
$ws.Cells.Item(55, 8).Value = 908.875
$ws.Cells.Item(55, 9).Value = 254.2
$ws.Cells.Item(55, 10).Value = 2000
$ws.Cells.Item(55, 11).Value = 762.5999999999999
$ws.Cells.Item(55, 12).Value = 6000
$ws.Cells.Item(55, 13).Value = -585.5999999999999
$ws.Cells.Item(55, 14).Value = -6354
$ws.Cells.Item(56, 8).Value = 10000
$ws.Cells.Item(56, 9).Value = 10000
$ws.Cells.Item(56, 11).Value = 10000
$ws.Cells.Item(56, 13).Value = -9470
$ws.Cells.Item(112, 8).Value = 6607.3335
$ws.Cells.Item(112, 9).Value = 6607.3335
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 19822.0005
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(112, 14).Value = -18714.0005
$ws.Cells.Item(131, 8).Value = 1374.75
$ws.Cells.Item(131, 10).Value = 1466.3334
$ws.Cells.Item(131, 12).Value = 4399.0002
$ws.Cells.Item(131, 14).Value = -14479.0002

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 22000
$ws.Cells.Item(33, 9).Value = 22000
$ws.Cells.Item(33, 11).Value = 22000
$ws.Cells.Item(33, 13).Value = -21748
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 12).ClearContents()
$ws.Cells.Item(44, 14).Value = 0
$ws.Cells.Item(47, 8).Value = 25500
$ws.Cells.Item(47, 10).Value = 25500
$ws.Cells.Item(47, 12).Value = 25500
$ws.Cells.Item(47, 14).Value = -26636
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(102, 8).Value = 1348.125
$ws.Cells.Item(102, 9).Value = 1162.5
$ws.Cells.Item(102, 10).Value = 1533.75
$ws.Cells.Item(102, 11).Value = 1162.5
$ws.Cells.Item(102, 12).Value = 1533.75
$ws.Cells.Item(102, 13).Value = 459.5
$ws.Cells.Item(102, 14).Value = -4777.75
$ws.Cells.Item(132, 8).Value = 2062.625
$ws.Cells.Item(132, 9).Value = 2062.625
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 6187.875
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -3657.875

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 8
$ws.Cells.Item(17, 9).Value = 8
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 8
$ws.Cells.Item(17, 12).ClearContents()
$ws.Cells.Item(17, 13).Value = 162
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(40, 8).Value = 4725.6
$ws.Cells.Item(40, 9).Value = 4608.909
$ws.Cells.Item(40, 11).Value = 4608.909
$ws.Cells.Item(40, 13).Value = -4472.909

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 2975
$ws.Cells.Item(4, 10).Value = 2975
$ws.Cells.Item(4, 12).Value = 2975
$ws.Cells.Item(4, 14).Value = -3201
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).ClearContents()
$ws.Cells.Item(74, 14).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).ClearContents()
$ws.Cells.Item(77, 14).Value = 0
$ws.Cells.Item(107, 8).Value = 441.76923
$ws.Cells.Item(107, 9).Value = 292.5
$ws.Cells.Item(107, 10).Value = 680.6
$ws.Cells.Item(107, 11).Value = 877.5
$ws.Cells.Item(107, 12).Value = 2041.8
$ws.Cells.Item(107, 13).Value = 1042.5
$ws.Cells.Item(107, 14).Value = -5881.8
$ws.Cells.Item(113, 8).Value = 347.875
$ws.Cells.Item(113, 9).Value = 356.33334
$ws.Cells.Item(113, 11).Value = 1069.00002
$ws.Cells.Item(113, 13).Value = 1100.99998
